$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Pruebas")

# --- Update "Pasos a ejectuar" (steps) text for TC-001 (H2) ---
$ws.Range("H2").Value = "1. Ingresar a www.saucedemo.com .`n2. Iniciar sesion con las credenciales validas.`n3. Verificar que el titulo de la pagina.`n4. Comprobar que exista al menos un producto visible.`n5. Validar la presencia de elementos importantes del interfaz (menu,filtro).`n6. Listar por nombre y/o precio."

# --- Update "Pasos a ejectuar" (steps) text for TC-002 (H3) ---
$ws.Range("H3").Value = "1. Ingresar a www.saucedemo.com .`n2. Iniciar sesion con las credenciales validas.`n3. Elegir un producto y hacer click en `"Add to cart`".`n4. Observar y hacer click en el icono del carrito .`n5. Verificar que el producto agregado se muestre en el carrito."

# --- Update "Observaciones" column (L2, L3): "Todo correcto." -> "Todo Ok" ---
$ws.Range("L2").Value = "Todo Ok"
$ws.Range("L3").Value = "Todo Ok"

# --- Update "Estado" column (K2, K3): "ok" -> "Ok" ---
$ws.Range("K2").Value = "Ok"
$ws.Range("K3").Value = "Ok"

# --- Normalize G2:G3 formatting (drop the stray empty-alignment style, reuse the
#     plain bordered style already used by A2:F2) ---
$ws.Range("A2").Copy()
$ws.Range("G2:G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Reset default (un-customized) row heights on the non-data rows ---
$ws.Rows("1").RowHeight = 12.75
$ws.Rows("4:20").RowHeight = 12.75

# --- Move active selection to A3 (matches saved sheetView selection) ---
$ws.Range("A3").Select()
